# Natmi following Dr Hou advice
# Rewrites the Dcn-Egfr LR-pair sheet to include the "ECs" cluster
# (4x4 cluster cross-product instead of 3x4) with updated NATMI stats.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$headers = @(
  "Sending cluster", "Ligand symbol", "Receptor symbol", "Target cluster", "Ligand-expressing cells", "Ligand detection rate", "Ligand average expression value", "Ligand total expression value", "Ligand derived specificity of average expression value", "Ligand derived specificity of total expression value", "Receptor-expressing cells", "Receptor detection rate", "Receptor average expression value", "Receptor total expression value", "Receptor derived specificity of average expression value", "Receptor derived specificity of total expression value", "Edge average expression weight", "Edge total expression weight", "Edge average expression derived specificity", "Edge total expression derived specificity"
)
for ($c = 1; $c -le $headers.Length; $c++) {
  $ws.Cells.Item(1, $c).Value = $headers[$c - 1]
}

# --- row 2 ---
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Dcn"
$ws.Cells.Item(2, 3).Value = "Egfr"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 2
$ws.Cells.Item(2, 6).Value = 0.6666666666666666
$ws.Cells.Item(2, 7).Value = 2.335690666666667
$ws.Cells.Item(2, 8).Value = 7.007072
$ws.Cells.Item(2, 9).Value = 0.0008591106367287776
$ws.Cells.Item(2, 10).Value = 0.0008591106367287777
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 2.701354
$ws.Cells.Item(2, 14).Value = 8.104061999999999
$ws.Cells.Item(2, 15).Value = 0.02221077311549548
$ws.Cells.Item(2, 16).Value = 0.02221077311549548
$ws.Cells.Item(2, 17).Value = 6.309527325162667
$ws.Cells.Item(2, 18).Value = 56.785745926464
$ws.Cells.Item(2, 19).Value = 0.00001908151143349174
$ws.Cells.Item(2, 20).Value = 0.00001908151143349174

# --- row 3 ---
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Dcn"
$ws.Cells.Item(3, 3).Value = "Egfr"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 2
$ws.Cells.Item(3, 6).Value = 0.6666666666666666
$ws.Cells.Item(3, 7).Value = 2.335690666666667
$ws.Cells.Item(3, 8).Value = 7.007072
$ws.Cells.Item(3, 9).Value = 0.0008591106367287776
$ws.Cells.Item(3, 10).Value = 0.0008591106367287777
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 88.14978533333333
$ws.Cells.Item(3, 14).Value = 264.449356
$ws.Cells.Item(3, 15).Value = 0.7247753838328104
$ws.Cells.Item(3, 16).Value = 0.7247753838328105
$ws.Cells.Item(3, 17).Value = 205.8906308717369
$ws.Cells.Item(3, 18).Value = 1853.015677845632
$ws.Cells.Item(3, 19).Value = 0.00062266224148995
$ws.Cells.Item(3, 20).Value = 0.0006226622414899501

# --- row 4 ---
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Dcn"
$ws.Cells.Item(4, 3).Value = "Egfr"
$ws.Cells.Item(4, 4).Value = "M2"
$ws.Cells.Item(4, 5).Value = 2
$ws.Cells.Item(4, 6).Value = 0.6666666666666666
$ws.Cells.Item(4, 7).Value = 2.335690666666667
$ws.Cells.Item(4, 8).Value = 7.007072
$ws.Cells.Item(4, 9).Value = 0.0008591106367287776
$ws.Cells.Item(4, 10).Value = 0.0008591106367287777
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 0.24063
$ws.Cells.Item(4, 14).Value = 0.72189
$ws.Cells.Item(4, 15).Value = 0.001978481285600361
$ws.Cells.Item(4, 16).Value = 0.001978481285600361
$ws.Cells.Item(4, 17).Value = 0.5620372451200001
$ws.Cells.Item(4, 18).Value = 5.05833520608
$ws.Cells.Item(4, 19).Value = 0.000001699734317028097
$ws.Cells.Item(4, 20).Value = 0.000001699734317028097

# --- row 5 ---
$ws.Cells.Item(5, 1).Value = "ECs"
$ws.Cells.Item(5, 2).Value = "Dcn"
$ws.Cells.Item(5, 3).Value = "Egfr"
$ws.Cells.Item(5, 4).Value = "sCs"
$ws.Cells.Item(5, 5).Value = 2
$ws.Cells.Item(5, 6).Value = 0.6666666666666666
$ws.Cells.Item(5, 7).Value = 2.335690666666667
$ws.Cells.Item(5, 8).Value = 7.007072
$ws.Cells.Item(5, 9).Value = 0.0008591106367287776
$ws.Cells.Item(5, 10).Value = 0.0008591106367287777
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 30.53182233333333
$ws.Cells.Item(5, 14).Value = 91.595467
$ws.Cells.Item(5, 15).Value = 0.2510353617660938
$ws.Cells.Item(5, 16).Value = 0.2510353617660938
$ws.Cells.Item(5, 17).Value = 71.31289246029156
$ws.Cells.Item(5, 18).Value = 641.816032142624
$ws.Cells.Item(5, 19).Value = 0.0002156671494883079
$ws.Cells.Item(5, 20).Value = 0.0002156671494883079

# --- row 6 ---
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Dcn"
$ws.Cells.Item(6, 3).Value = "Egfr"
$ws.Cells.Item(6, 4).Value = "ECs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 2706.934895666667
$ws.Cells.Item(6, 8).Value = 8120.804687
$ws.Cells.Item(6, 9).Value = 0.9956611956318718
$ws.Cells.Item(6, 10).Value = 0.9956611956318719
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 2.701354
$ws.Cells.Item(6, 14).Value = 8.104061999999999
$ws.Cells.Item(6, 15).Value = 0.02221077311549548
$ws.Cells.Item(6, 16).Value = 0.02221077311549548
$ws.Cells.Item(6, 17).Value = 7312.389408148732
$ws.Cells.Item(6, 18).Value = 65811.50467333858
$ws.Cells.Item(6, 19).Value = 0.02211440491608246
$ws.Cells.Item(6, 20).Value = 0.02211440491608247

# --- row 7 ---
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Dcn"
$ws.Cells.Item(7, 3).Value = "Egfr"
$ws.Cells.Item(7, 4).Value = "FAPs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 2706.934895666667
$ws.Cells.Item(7, 8).Value = 8120.804687
$ws.Cells.Item(7, 9).Value = 0.9956611956318718
$ws.Cells.Item(7, 10).Value = 0.9956611956318719
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 88.14978533333333
$ws.Cells.Item(7, 14).Value = 264.449356
$ws.Cells.Item(7, 15).Value = 0.7247753838328104
$ws.Cells.Item(7, 16).Value = 0.7247753838328105
$ws.Cells.Item(7, 17).Value = 238615.7299643257
$ws.Cells.Item(7, 18).Value = 2147541.569678931
$ws.Cells.Item(7, 19).Value = 0.7216307252315248
$ws.Cells.Item(7, 20).Value = 0.721630725231525

# --- row 8 ---
$ws.Cells.Item(8, 1).Value = "FAPs"
$ws.Cells.Item(8, 2).Value = "Dcn"
$ws.Cells.Item(8, 3).Value = "Egfr"
$ws.Cells.Item(8, 4).Value = "M2"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 2706.934895666667
$ws.Cells.Item(8, 8).Value = 8120.804687
$ws.Cells.Item(8, 9).Value = 0.9956611956318718
$ws.Cells.Item(8, 10).Value = 0.9956611956318719
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 0.24063
$ws.Cells.Item(8, 14).Value = 0.72189
$ws.Cells.Item(8, 15).Value = 0.001978481285600361
$ws.Cells.Item(8, 16).Value = 0.001978481285600361
$ws.Cells.Item(8, 17).Value = 651.36974394427
$ws.Cells.Item(8, 18).Value = 5862.32769549843
$ws.Cells.Item(8, 19).Value = 0.001969897042356139
$ws.Cells.Item(8, 20).Value = 0.001969897042356139

# --- row 9 ---
$ws.Cells.Item(9, 1).Value = "FAPs"
$ws.Cells.Item(9, 2).Value = "Dcn"
$ws.Cells.Item(9, 3).Value = "Egfr"
$ws.Cells.Item(9, 4).Value = "sCs"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 2706.934895666667
$ws.Cells.Item(9, 8).Value = 8120.804687
$ws.Cells.Item(9, 9).Value = 0.9956611956318718
$ws.Cells.Item(9, 10).Value = 0.9956611956318719
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 30.53182233333333
$ws.Cells.Item(9, 14).Value = 91.595467
$ws.Cells.Item(9, 15).Value = 0.2510353617660938
$ws.Cells.Item(9, 16).Value = 0.2510353617660938
$ws.Cells.Item(9, 17).Value = 82647.65530239488
$ws.Cells.Item(9, 18).Value = 743828.8977215538
$ws.Cells.Item(9, 19).Value = 0.2499461684419085
$ws.Cells.Item(9, 20).Value = 0.2499461684419085

# --- row 10 ---
$ws.Cells.Item(10, 1).Value = "M2"
$ws.Cells.Item(10, 2).Value = "Dcn"
$ws.Cells.Item(10, 3).Value = "Egfr"
$ws.Cells.Item(10, 4).Value = "ECs"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 1.085526333333333
$ws.Cells.Item(10, 8).Value = 3.256579
$ws.Cells.Item(10, 9).Value = 0.0003992768531916849
$ws.Cells.Item(10, 10).Value = 0.0003992768531916849
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 2.701354
$ws.Cells.Item(10, 14).Value = 8.104061999999999
$ws.Cells.Item(10, 15).Value = 0.02221077311549548
$ws.Cells.Item(10, 16).Value = 0.02221077311549548
$ws.Cells.Item(10, 17).Value = 2.932390902655333
$ws.Cells.Item(10, 18).Value = 26.391518123898
$ws.Cells.Item(10, 19).Value = 0.00000886824759650951
$ws.Cells.Item(10, 20).Value = 0.000008868247596509512

# --- row 11 ---
$ws.Cells.Item(11, 1).Value = "M2"
$ws.Cells.Item(11, 2).Value = "Dcn"
$ws.Cells.Item(11, 3).Value = "Egfr"
$ws.Cells.Item(11, 4).Value = "FAPs"
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(11, 7).Value = 1.085526333333333
$ws.Cells.Item(11, 8).Value = 3.256579
$ws.Cells.Item(11, 9).Value = 0.0003992768531916849
$ws.Cells.Item(11, 10).Value = 0.0003992768531916849
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 12).Value = 1
$ws.Cells.Item(11, 13).Value = 88.14978533333333
$ws.Cells.Item(11, 14).Value = 264.449356
$ws.Cells.Item(11, 15).Value = 0.7247753838328104
$ws.Cells.Item(11, 16).Value = 0.7247753838328105
$ws.Cells.Item(11, 17).Value = 95.68891325701378
$ws.Cells.Item(11, 18).Value = 861.200219313124
$ws.Cells.Item(11, 19).Value = 0.0002893860345275601
$ws.Cells.Item(11, 20).Value = 0.0002893860345275602

# --- row 12 ---
$ws.Cells.Item(12, 1).Value = "M2"
$ws.Cells.Item(12, 2).Value = "Dcn"
$ws.Cells.Item(12, 3).Value = "Egfr"
$ws.Cells.Item(12, 4).Value = "M2"
$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 6).Value = 1
$ws.Cells.Item(12, 7).Value = 1.085526333333333
$ws.Cells.Item(12, 8).Value = 3.256579
$ws.Cells.Item(12, 9).Value = 0.0003992768531916849
$ws.Cells.Item(12, 10).Value = 0.0003992768531916849
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 12).Value = 1
$ws.Cells.Item(12, 13).Value = 0.24063
$ws.Cells.Item(12, 14).Value = 0.72189
$ws.Cells.Item(12, 15).Value = 0.001978481285600361
$ws.Cells.Item(12, 16).Value = 0.001978481285600361
$ws.Cells.Item(12, 17).Value = 0.26121020159
$ws.Cells.Item(12, 18).Value = 2.35089181431
$ws.Cells.Item(12, 19).Value = 0.0000007899617818131514
$ws.Cells.Item(12, 20).Value = 0.0000007899617818131515

# --- row 13 ---
$ws.Cells.Item(13, 1).Value = "M2"
$ws.Cells.Item(13, 2).Value = "Dcn"
$ws.Cells.Item(13, 3).Value = "Egfr"
$ws.Cells.Item(13, 4).Value = "sCs"
$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 6).Value = 1
$ws.Cells.Item(13, 7).Value = 1.085526333333333
$ws.Cells.Item(13, 8).Value = 3.256579
$ws.Cells.Item(13, 9).Value = 0.0003992768531916849
$ws.Cells.Item(13, 10).Value = 0.0003992768531916849
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 12).Value = 1
$ws.Cells.Item(13, 13).Value = 30.53182233333333
$ws.Cells.Item(13, 14).Value = 91.595467
$ws.Cells.Item(13, 15).Value = 0.2510353617660938
$ws.Cells.Item(13, 16).Value = 0.2510353617660938
$ws.Cells.Item(13, 17).Value = 33.14309714748811
$ws.Cells.Item(13, 18).Value = 298.287874327393
$ws.Cells.Item(13, 19).Value = 0.0001002326092858021
$ws.Cells.Item(13, 20).Value = 0.0001002326092858022

# --- row 14 ---
$ws.Cells.Item(14, 1).Value = "sCs"
$ws.Cells.Item(14, 2).Value = "Dcn"
$ws.Cells.Item(14, 3).Value = "Egfr"
$ws.Cells.Item(14, 4).Value = "ECs"
$ws.Cells.Item(14, 5).Value = 3
$ws.Cells.Item(14, 6).Value = 1
$ws.Cells.Item(14, 7).Value = 8.374824666666667
$ws.Cells.Item(14, 8).Value = 25.124474
$ws.Cells.Item(14, 9).Value = 0.003080416878207562
$ws.Cells.Item(14, 10).Value = 0.003080416878207562
$ws.Cells.Item(14, 11).Value = 3
$ws.Cells.Item(14, 12).Value = 1
$ws.Cells.Item(14, 13).Value = 2.701354
$ws.Cells.Item(14, 14).Value = 8.104061999999999
$ws.Cells.Item(14, 15).Value = 0.02221077311549548
$ws.Cells.Item(14, 16).Value = 0.02221077311549548
$ws.Cells.Item(14, 17).Value = 22.62336611259867
$ws.Cells.Item(14, 18).Value = 203.610295013388
$ws.Cells.Item(14, 19).Value = 0.00006841844038301103
$ws.Cells.Item(14, 20).Value = 0.00006841844038301103

# --- row 15 ---
$ws.Cells.Item(15, 1).Value = "sCs"
$ws.Cells.Item(15, 2).Value = "Dcn"
$ws.Cells.Item(15, 3).Value = "Egfr"
$ws.Cells.Item(15, 4).Value = "FAPs"
$ws.Cells.Item(15, 5).Value = 3
$ws.Cells.Item(15, 6).Value = 1
$ws.Cells.Item(15, 7).Value = 8.374824666666667
$ws.Cells.Item(15, 8).Value = 25.124474
$ws.Cells.Item(15, 9).Value = 0.003080416878207562
$ws.Cells.Item(15, 10).Value = 0.003080416878207562
$ws.Cells.Item(15, 11).Value = 3
$ws.Cells.Item(15, 12).Value = 1
$ws.Cells.Item(15, 13).Value = 88.14978533333333
$ws.Cells.Item(15, 14).Value = 264.449356
$ws.Cells.Item(15, 15).Value = 0.7247753838328104
$ws.Cells.Item(15, 16).Value = 0.7247753838328105
$ws.Cells.Item(15, 17).Value = 738.2389965709715
$ws.Cells.Item(15, 18).Value = 6644.150969138743
$ws.Cells.Item(15, 19).Value = 0.002232610325267954
$ws.Cells.Item(15, 20).Value = 0.002232610325267954

# --- row 16 ---
$ws.Cells.Item(16, 1).Value = "sCs"
$ws.Cells.Item(16, 2).Value = "Dcn"
$ws.Cells.Item(16, 3).Value = "Egfr"
$ws.Cells.Item(16, 4).Value = "M2"
$ws.Cells.Item(16, 5).Value = 3
$ws.Cells.Item(16, 6).Value = 1
$ws.Cells.Item(16, 7).Value = 8.374824666666667
$ws.Cells.Item(16, 8).Value = 25.124474
$ws.Cells.Item(16, 9).Value = 0.003080416878207562
$ws.Cells.Item(16, 10).Value = 0.003080416878207562
$ws.Cells.Item(16, 11).Value = 3
$ws.Cells.Item(16, 12).Value = 1
$ws.Cells.Item(16, 13).Value = 0.24063
$ws.Cells.Item(16, 14).Value = 0.72189
$ws.Cells.Item(16, 15).Value = 0.001978481285600361
$ws.Cells.Item(16, 16).Value = 0.001978481285600361
$ws.Cells.Item(16, 17).Value = 2.01523405954
$ws.Cells.Item(16, 18).Value = 18.13710653586
$ws.Cells.Item(16, 19).Value = 0.000006094547145381149
$ws.Cells.Item(16, 20).Value = 0.000006094547145381149

# --- row 17 ---
$ws.Cells.Item(17, 1).Value = "sCs"
$ws.Cells.Item(17, 2).Value = "Dcn"
$ws.Cells.Item(17, 3).Value = "Egfr"
$ws.Cells.Item(17, 4).Value = "sCs"
$ws.Cells.Item(17, 5).Value = 3
$ws.Cells.Item(17, 6).Value = 1
$ws.Cells.Item(17, 7).Value = 8.374824666666667
$ws.Cells.Item(17, 8).Value = 25.124474
$ws.Cells.Item(17, 9).Value = 0.003080416878207562
$ws.Cells.Item(17, 10).Value = 0.003080416878207562
$ws.Cells.Item(17, 11).Value = 3
$ws.Cells.Item(17, 12).Value = 1
$ws.Cells.Item(17, 13).Value = 30.53182233333333
$ws.Cells.Item(17, 14).Value = 91.595467
$ws.Cells.Item(17, 15).Value = 0.2510353617660938
$ws.Cells.Item(17, 16).Value = 0.2510353617660938
$ws.Cells.Item(17, 17).Value = 255.6986587954842
$ws.Cells.Item(17, 18).Value = 2301.287929159358
$ws.Cells.Item(17, 19).Value = 0.0007732935654112167
$ws.Cells.Item(17, 20).Value = 0.0007732935654112167

